# "Worked on temporal resolution" -- extend the Demand sheet (EU27.Elec)
# to the full set of time steps (t = 0..12), matching the resolution
# already used on the SupIm sheet, and make the Demand sheet the active
# / selected sheet (previously SupIm was the selected tab).

$wb = $excel.ActiveWorkbook

$demand = $wb.Worksheets.Item("Demand")

# --- Extend data ---------------------------------------------------------
# Copy the number formatting of the existing data row (row 3) down to the
# new rows so the newly added cells pick up the same style (thousands
# separator) as the rest of the column.
$demand.Range("A3:B3").Copy()
$demand.Range("A4:B14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update existing value and fill in the new time steps / values.
$demand.Range("B3").Value = 508833333

$demand.Range("A4").Value = 2
$demand.Range("A5").Value = 3
$demand.Range("A6").Value = 4
$demand.Range("A7").Value = 5
$demand.Range("A8").Value = 6
$demand.Range("A9").Value = 7
$demand.Range("A10").Value = 8
$demand.Range("A11").Value = 9
$demand.Range("A12").Value = 10
$demand.Range("A13").Value = 11
$demand.Range("A14").Value = 12

$demand.Range("B4").Value = 508833333
$demand.Range("B5").Value = 508833333
$demand.Range("B6").Value = 508833333
$demand.Range("B7").Value = 508833333
$demand.Range("B8").Value = 508833333
$demand.Range("B9").Value = 508833333
$demand.Range("B10").Value = 508833333
$demand.Range("B11").Value = 508833333
$demand.Range("B12").Value = 508833333
$demand.Range("B13").Value = 508833333
$demand.Range("B14").Value = 508833333

# Auto fit column B now that it holds longer numbers.
$demand.Range("B1:B14").EntireColumn.AutoFit()

# --- Selection / active sheet ---------------------------------------------
# Previously SupIm was the selected/active tab; now Demand is.
$demand.Activate()
$demand.Range("D14").Select()
